# Update "想去人数" (wanted-to-go count) figures on both the "展览" sheet
# and the aggregated "全部类型" sheet, matching a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (index 1) ---
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F3").Value  = 50
$wsExhibit.Range("F4").Value  = 1448
$wsExhibit.Range("F5").Value  = 334
$wsExhibit.Range("F6").Value  = 1049
$wsExhibit.Range("F7").Value  = 10829
$wsExhibit.Range("F12").Value = 729
$wsExhibit.Range("F13").Value = 12143
$wsExhibit.Range("F14").Value = 12608

# --- Sheet 4: 全部类型 (index 4) ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value  = 50
$wsAll.Range("F5").Value  = 1448
$wsAll.Range("F6").Value  = 334
$wsAll.Range("F7").Value  = 1049
$wsAll.Range("F8").Value  = 10829
$wsAll.Range("F13").Value = 729
$wsAll.Range("F14").Value = 12143
$wsAll.Range("F15").Value = 12608
